# Auto-generated edit script: updates Betfair Back/Lay odds cells
# to match the target snapshot for 2025-12-26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.95
$ws.Range("J2").Value = 4.4
$ws.Range("K2").Value = 4.6
$ws.Range("P2").Value = 2.9
$ws.Range("G3").Value = 1.96
$ws.Range("H3").Value = 1.09
$ws.Range("J3").Value = 2.08
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 1.65
$ws.Range("H4").Value = 1.63
$ws.Range("I4").Value = 1.76
$ws.Range("J4").Value = 4.3
$ws.Range("N4").Value = 5.4
$ws.Range("O4").Value = 1.17
$ws.Range("G5").Value = 2.46
$ws.Range("H5").Value = 3.85
$ws.Range("K5").Value = 3.55
$ws.Range("G6").Value = 1.5
$ws.Range("J6").Value = 3.95
$ws.Range("K6").Value = 4.9
$ws.Range("N6").Value = 3
$ws.Range("W6").Value = 2.96
$ws.Range("G7").Value = 1.13
$ws.Range("W7").Value = 7.6
$ws.Range("I8").Value = 12
$ws.Range("R8").Value = 1.75
$ws.Range("F9").Value = 3.2
$ws.Range("H9").Value = 2.16
$ws.Range("I9").Value = 2.4
$ws.Range("P9").Value = 1.97
$ws.Range("Q9").Value = 1.78
$ws.Range("V9").Value = 1.71
$ws.Range("G12").Value = 2.64
$ws.Range("W12").Value = 1.61
$ws.Range("X12").Value = 14.5
$ws.Range("Y12").Value = 14
$ws.Range("AB12").Value = 11.5
$ws.Range("AC12").Value = 9.199999999999999
$ws.Range("AG12").Value = 14.5
$ws.Range("H13").Value = 1.46
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 5.3
$ws.Range("R13").Value = 1.58
$ws.Range("S13").Value = 2.38
$ws.Range("W13").Value = 1.13
$ws.Range("X13").Value = 32
$ws.Range("Y13").Value = 13.5
$ws.Range("Z13").Value = 12.5
$ws.Range("AA13").Value = 16.5
$ws.Range("AB13").Value = 38
$ws.Range("AC13").Value = 14.5
$ws.Range("AD13").Value = 13
$ws.Range("AE13").Value = 17.5
$ws.Range("AG13").Value = 34
$ws.Range("AH13").Value = 27
$ws.Range("AI13").Value = 38
$ws.Range("AO13").Value = 6.8
$ws.Range("K14").Value = 4.1
$ws.Range("R15").Value = 1.85
$ws.Range("S15").Value = 1.94
$ws.Range("T15").Value = 1.8
$ws.Range("U15").Value = 1.98
$ws.Range("Z15").Value = 150
$ws.Range("AE15").Value = 190
$ws.Range("AF15").Value = 11.5
$ws.Range("AI15").Value = 140
$ws.Range("AJ15").Value = 11.5
$ws.Range("AM15").Value = 140
$ws.Range("AN15").Value = 3.65
$ws.Range("AO15").Value = 190
$ws.Range("F16").Value = 2.1
$ws.Range("G16").Value = 2.34
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 3.9
$ws.Range("J16").Value = 3.75
$ws.Range("O16").Value = 1.26
$ws.Range("Q16").Value = 1.76
$ws.Range("T16").Value = 1.65
$ws.Range("V16").Value = 1.37
$ws.Range("W16").Value = 1.75
$ws.Range("X16").Value = 21
$ws.Range("AA16").Value = 65
$ws.Range("O18").Value = 1.3
$ws.Range("AK18").Value = 38
$ws.Range("T19").Value = 1.51
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 1.48
$ws.Range("Q20").Value = 1.38
$ws.Range("R20").Value = 1.83
$ws.Range("S20").Value = 1.99
$ws.Range("V20").Value = 2.7
$ws.Range("O21").Value = 1.25
$ws.Range("G22").Value = 2.58
$ws.Range("H22").Value = 3.45
$ws.Range("I22").Value = 3.75
$ws.Range("Q22").Value = 2.48
$ws.Range("K23").Value = 5.4
$ws.Range("N23").Value = 2.4
$ws.Range("P23").Value = 2.4
$ws.Range("Q23").Value = 1.57
$ws.Range("K24").Value = 9.199999999999999
$ws.Range("Q24").Value = 1.3
$ws.Range("R24").Value = 2.14
$ws.Range("F25").Value = 2.66
$ws.Range("G25").Value = 2.68
$ws.Range("O25").Value = 1.22
$ws.Range("Q25").Value = 1.68
$ws.Range("W25").Value = 1.59
$ws.Range("AI25").Value = 32
